$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns for the rows we touch,
# so Excel does not auto-coerce values like "1.00" or "0.100" into numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

# Apply updated values
$ws.Range('D2').Value = '36.629.64'
$ws.Range('D3').Value = '1.963.68'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '244.08'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').Value = '58.72'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.375'
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('D10').Value = '0.0802'
$ws.Range('E10').Value = '  -4.05%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  +3.33%  '
$ws.Range('D13').Value = '2.253.50'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '0.822'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '13.66'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').Value = '1.969.11'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').Value = '36.577.80'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '69.67'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '0.0₃0859'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '228.72'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.07'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('D26').Value = '9.31'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '0.139'
$ws.Range('E27').Value = '  +11.43%  '
$ws.Range('D28').Value = '160.29'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').Value = '19.36'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('D32').Value = '4.68'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = '6.09'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.24'
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '3.39'
$ws.Range('E38').Value = '  +13.25%  '
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  +3.39%  '
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0212'
$ws.Range('E42').Value = '  +1.21%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.16'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('D44').Value = '16.02'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = '1.359.63'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').Value = '87.51'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = '7.13'
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').Value = '2.144.54'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  -5.54%  '

Write-Output "Applied crypto price/volume updates"
